$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for data rows 2-21 (row 1 is the untouched header).
# This mirrors the diff: 3 brand-new rows are inserted at the top (new rows
# 2-4), the previous rows 2-17 shift down to rows 5-20, the old row 18 is
# duplicated as the new last row (21), and the old rows 19-20 are dropped.
$targetValues = @(
    @(0.0371100641787052, -0.0024434609804302, 0.0224492978304624),
    @(0.0600175112485885, -0.0474947728216648, 0.0047342055477201),
    @(0.0438295826315879, -0.08277224004268641, 0.07590000331401819),
    @(-0.06276640295982359, -0.5068654417991638, -0.0064140851609408),
    @(-0.0687223374843597, -0.2246456891298294, 0.1652390509843826),
    @(-0.2229658216238021, -0.113315500319004, 0.2180788964033126),
    @(0.3182607889175415, 0.6488916277885437, 0.4970915913581848),
    @(0.3697261810302734, 1.551903128623962, -0.230448916554451),
    @(0.451123982667923, 0.8857545852661133, -1.012661814689636),
    @(-0.1492038369178772, -1.078635334968567, -4.738939762115479),
    @(-1.247386813163757, -1.642922043800354, 0.3756821155548095),
    @(2.319607973098755, -0.2393064647912979, -2.242333650588989),
    @(0.1510364264249801, 2.41413950920105, 0.4211915731430053),
    @(0.2503020465373993, -1.337489485740662, 0.0983493030071258),
    @(-0.4132503271102905, 0.2756529450416565, -0.5613851547241211),
    @(-0.1337794959545135, 0.3077233731746673, -0.0181732401251792),
    @(0.0068722339347004, -0.2590068578720093, 0.2180788964033126),
    @(-0.0940732508897781, -0.0920879393815994, 0.011148290708661),
    @(-0.0429132841527462, -0.0099265603348612, -0.2417499274015426),
    @(0.0386372283101081, 0.07590000331401819, 0.0108428578823804)
)

for ($i = 0; $i -lt $targetValues.Length; $i++) {
    $r = 2 + $i
    $values = $targetValues[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
}
